# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" between "总计" and "2021-Q4",
#    populated with that quarter's fund-holding data. It is created by
#    duplicating the existing "2021-Q4" sheet (so it inherits the same
#    look: bold/centered/bordered header row + index column) and then
#    overwriting its cell values.
# 2) Update the "总计" (summary) sheet: the new 2022-Q3 totals become row
#    2, and the pre-existing 2021-Q4 totals shift down to row 3 (with its
#    index column bumped from 0 to 1).

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

# --- Step 1: add the new "2022-Q3" sheet right after "总计" ------------
# Duplicate "2021-Q4" (placed right after "总计") so the new sheet starts
# out with identical formatting/styles, then rename + overwrite its data.
$q4.Copy([System.Reflection.Missing]::Value, $total)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# --- helper pattern: write a text cell without Excel's "looks like a
# number" auto-coercion turning e.g. "011160" into 11160, or "3.70" into
# 3.7. Force a text number format for the assignment, then clear the
# format again so no stray style is left on the cell. ---

# Row 2
$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "011160"
$q3.Range("B2").ClearFormats()
$q3.Range("C2").NumberFormat = "@"
$q3.Range("C2").Value = "富国质量成长6个月持有期混合A"
$q3.Range("C2").ClearFormats()
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "3.70"
$q3.Range("D2").ClearFormats()
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "85.89"
$q3.Range("E2").ClearFormats()
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "3.88"
$q3.Range("F2").ClearFormats()
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.1436"
$q3.Range("G2").ClearFormats()
$q3.Range("H2").Value = 4

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "004448"
$q3.Range("B3").ClearFormats()
$q3.Range("C3").NumberFormat = "@"
$q3.Range("C3").Value = "博时汇智回报灵活配置混合"
$q3.Range("C3").ClearFormats()
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "1.77"
$q3.Range("D3").ClearFormats()
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "67.69"
$q3.Range("E3").ClearFormats()
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "3.41"
$q3.Range("F3").ClearFormats()
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0604"
$q3.Range("G3").ClearFormats()
$q3.Range("H3").Value = 6

# Row 4 (new row - copy column-A styling down from row 3 first so it
# matches the bold/centered/bordered index-column look)
$q3.Range("A3").Copy()
$q3.Range("A4").PasteSpecial(-4122)
$q3.Range("A4").Value = 2
$q3.Range("B4").NumberFormat = "@"
$q3.Range("B4").Value = "011161"
$q3.Range("B4").ClearFormats()
$q3.Range("C4").NumberFormat = "@"
$q3.Range("C4").Value = "富国质量成长6个月持有期混合C"
$q3.Range("C4").ClearFormats()
$q3.Range("D4").NumberFormat = "@"
$q3.Range("D4").Value = "0.14"
$q3.Range("D4").ClearFormats()
$q3.Range("E4").NumberFormat = "@"
$q3.Range("E4").Value = "85.89"
$q3.Range("E4").ClearFormats()
$q3.Range("F4").NumberFormat = "@"
$q3.Range("F4").Value = "3.88"
$q3.Range("F4").ClearFormats()
$q3.Range("G4").NumberFormat = "@"
$q3.Range("G4").Value = "0.0054"
$q3.Range("G4").ClearFormats()
$q3.Range("H4").Value = 4

# --- Step 2: update "总计" - shift 2021-Q4 row down, add 2022-Q3 row --
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.1
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.21

Write-Host "applied 2022-Q3 edit"
